# Update column F (dSF) values per the commit "repull data, push all data, mean calculation"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 4
    3  = -1
    4  = 3
    6  = -2
    7  = 6
    8  = 1
    9  = -2
    10 = -2
    11 = 0
    13 = 3
    14 = 3
    15 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
